# Apply the meta property change: orientation=portrait -> orientation=c
# (excel meta "orientation" no longer supports landscape|l|portrait|p)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the B1 cell that stores the orientation meta value.
$ws.Range("B1").Value = "orientation=c"

# Move the selection to B2 (was B13 before the edit).
$ws.Range("B2").Select()
